# edit.ps1
# Applies the re-scrape update to the Bundesliga 2023-2024 odds sheet:
#  1) Re-orders the match-detail columns (F:V) for several same-date groups
#     of rows so they reflect the newly re-scraped ordering (index/date
#     columns A:E for each row stay put - only F:V travel).
#  2) Appends two brand-new matches as rows 90 and 91 (Indice 89 and 90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Reorder match rows within same-date groups to match the re-scraped order ---
# Row 23
$ws.Range("F23").Value = "Bayer Leverkusen"
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = "Darmstadt"
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 1.4
$ws.Range("K23").Value = "21/08/2023 06:14"
$ws.Range("L23").Value = 1.22
$ws.Range("M23").Value = "02/09/2023 15:24"
$ws.Range("N23").Value = 4.92
$ws.Range("O23").Value = "21/08/2023 06:14"
$ws.Range("P23").Value = 7.25
$ws.Range("Q23").Value = "02/09/2023 15:24"
$ws.Range("R23").Value = 6.76
$ws.Range("S23").Value = "21/08/2023 06:14"
$ws.Range("T23").Value = 13
$ws.Range("U23").Value = "02/09/2023 15:24"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayer-leverkusen-darmstadt/8tjtyQa9/"

# Row 24
$ws.Range("F24").Value = "Augsburg"
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = "Bochum"
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 1.92
$ws.Range("K24").Value = "20/08/2023 09:02"
$ws.Range("L24").Value = 2.25
$ws.Range("M24").Value = "02/09/2023 15:16"
$ws.Range("N24").Value = 3.74
$ws.Range("O24").Value = "20/08/2023 09:02"
$ws.Range("P24").Value = 3.68
$ws.Range("Q24").Value = "02/09/2023 15:16"
$ws.Range("R24").Value = 3.7
$ws.Range("S24").Value = "20/08/2023 09:02"
$ws.Range("T24").Value = 3.26
$ws.Range("U24").Value = "02/09/2023 15:16"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-bochum/vVTOS17k/"

# Row 30
$ws.Range("F30").Value = "RB Leipzig"
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = "Augsburg"
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1.2
$ws.Range("K30").Value = "28/08/2023 16:01"
$ws.Range("L30").Value = 1.27
$ws.Range("M30").Value = "16/09/2023 15:25"
$ws.Range("N30").Value = 7.5
$ws.Range("O30").Value = "28/08/2023 16:01"
$ws.Range("P30").Value = 6.75
$ws.Range("Q30").Value = "16/09/2023 15:19"
$ws.Range("R30").Value = 13.58
$ws.Range("S30").Value = "28/08/2023 16:01"
$ws.Range("T30").Value = 10.25
$ws.Range("U30").Value = "16/09/2023 15:25"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/germany/bundesliga/rb-leipzig-augsburg/l0BFkPE2/"

# Row 31
$ws.Range("F31").Value = "FC Koln"
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = "Hoffenheim"
$ws.Range("I31").Value = 3
$ws.Range("J31").Value = 2.03
$ws.Range("K31").Value = "28/08/2023 16:01"
$ws.Range("L31").Value = 2.22
$ws.Range("M31").Value = "16/09/2023 15:17"
$ws.Range("N31").Value = 3.7
$ws.Range("O31").Value = "28/08/2023 16:01"
$ws.Range("P31").Value = 3.86
$ws.Range("Q31").Value = "16/09/2023 15:27"
$ws.Range("R31").Value = 3.68
$ws.Range("S31").Value = "28/08/2023 16:01"
$ws.Range("T31").Value = 3.19
$ws.Range("U31").Value = "16/09/2023 15:20"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-hoffenheim/lbWeVfTm/"

# Row 32
$ws.Range("F32").Value = "Freiburg"
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = "Dortmund"
$ws.Range("I32").Value = 4
$ws.Range("J32").Value = 3.62
$ws.Range("K32").Value = "28/08/2023 16:01"
$ws.Range("L32").Value = 2.97
$ws.Range("M32").Value = "16/09/2023 15:29"
$ws.Range("N32").Value = 3.92
$ws.Range("O32").Value = "28/08/2023 16:01"
$ws.Range("P32").Value = 3.82
$ws.Range("Q32").Value = "16/09/2023 15:27"
$ws.Range("R32").Value = 1.9
$ws.Range("S32").Value = "28/08/2023 16:01"
$ws.Range("T32").Value = 2.36
$ws.Range("U32").Value = "16/09/2023 15:28"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-dortmund/0zFJlqU8/"

# Row 33
$ws.Range("F33").Value = "Mainz"
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = "Stuttgart"
$ws.Range("I33").Value = 3
$ws.Range("J33").Value = 2.62
$ws.Range("K33").Value = "28/08/2023 16:01"
$ws.Range("L33").Value = 2.71
$ws.Range("M33").Value = "16/09/2023 15:25"
$ws.Range("N33").Value = 3.4
$ws.Range("O33").Value = "28/08/2023 16:01"
$ws.Range("P33").Value = 3.54
$ws.Range("Q33").Value = "16/09/2023 15:27"
$ws.Range("R33").Value = 2.79
$ws.Range("S33").Value = "28/08/2023 16:01"
$ws.Range("T33").Value = 2.71
$ws.Range("U33").Value = "16/09/2023 15:27"
$ws.Range("V33").Value = "https://www.betexplorer.com/football/germany/bundesliga/mainz-vfb-stuttgart/WMyiWzEs/"

# Row 34
$ws.Range("F34").Value = "Wolfsburg"
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = "Union Berlin"
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = 2.19
$ws.Range("K34").Value = "28/08/2023 16:01"
$ws.Range("L34").Value = 2.23
$ws.Range("M34").Value = "16/09/2023 15:27"
$ws.Range("N34").Value = 3.35
$ws.Range("O34").Value = "28/08/2023 16:01"
$ws.Range("P34").Value = 3.35
$ws.Range("Q34").Value = "16/09/2023 15:27"
$ws.Range("R34").Value = 3.33
$ws.Range("S34").Value = "28/08/2023 16:01"
$ws.Range("T34").Value = 3.63
$ws.Range("U34").Value = "16/09/2023 15:29"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/germany/bundesliga/wolfsburg-union-berlin/fcENm3qF/"

# Row 39
$ws.Range("F39").Value = "Augsburg"
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = "Mainz"
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 2.6
$ws.Range("K39").Value = "05/09/2023 12:01"
$ws.Range("L39").Value = 2.56
$ws.Range("M39").Value = "23/09/2023 15:25"
$ws.Range("N39").Value = 3.57
$ws.Range("O39").Value = "05/09/2023 12:01"
$ws.Range("P39").Value = 3.52
$ws.Range("Q39").Value = "23/09/2023 15:05"
$ws.Range("R39").Value = 2.81
$ws.Range("S39").Value = "05/09/2023 12:01"
$ws.Range("T39").Value = 2.88
$ws.Range("U39").Value = "23/09/2023 15:25"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-mainz/tz4tshSa/"

# Row 40
$ws.Range("F40").Value = "Bayern Munich"
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = "Bochum"
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1.07
$ws.Range("K40").Value = "05/09/2023 12:01"
$ws.Range("L40").Value = 1.12
$ws.Range("M40").Value = "23/09/2023 15:00"
$ws.Range("N40").Value = 20.79
$ws.Range("O40").Value = "05/09/2023 12:01"
$ws.Range("P40").Value = 10.5
$ws.Range("Q40").Value = "23/09/2023 14:59"
$ws.Range("R40").Value = 29.26
$ws.Range("S40").Value = "05/09/2023 12:01"
$ws.Range("T40").Value = 21
$ws.Range("U40").Value = "23/09/2023 15:00"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayern-munich-bochum/IwOCRCSC/"

# Row 41
$ws.Range("F41").Value = "Union Berlin"
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = "Hoffenheim"
$ws.Range("I41").Value = 2
$ws.Range("J41").Value = 1.85
$ws.Range("K41").Value = "05/09/2023 12:01"
$ws.Range("L41").Value = 2.13
$ws.Range("M41").Value = "23/09/2023 15:05"
$ws.Range("N41").Value = 3.93
$ws.Range("O41").Value = "05/09/2023 12:01"
$ws.Range("P41").Value = 3.72
$ws.Range("Q41").Value = "23/09/2023 15:20"
$ws.Range("R41").Value = 4.38
$ws.Range("S41").Value = "05/09/2023 12:01"
$ws.Range("T41").Value = 3.51
$ws.Range("U41").Value = "23/09/2023 15:26"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/germany/bundesliga/union-berlin-hoffenheim/AoQckGKJ/"

# Row 42
$ws.Range("F42").Value = "Dortmund"
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = "Wolfsburg"
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 1.47
$ws.Range("K42").Value = "05/09/2023 12:01"
$ws.Range("L42").Value = 1.74
$ws.Range("M42").Value = "23/09/2023 15:20"
$ws.Range("N42").Value = 5.02
$ws.Range("O42").Value = "05/09/2023 12:01"
$ws.Range("P42").Value = 4.42
$ws.Range("Q42").Value = "23/09/2023 15:28"
$ws.Range("R42").Value = 5.4
$ws.Range("S42").Value = "05/09/2023 12:01"
$ws.Range("T42").Value = 4.48
$ws.Range("U42").Value = "23/09/2023 15:28"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/germany/bundesliga/dortmund-wolfsburg/dUMGQWsJ/"

# Row 43
$ws.Range("F43").Value = "B. Monchengladbach"
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = "RB Leipzig"
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = 4.05
$ws.Range("K43").Value = "05/09/2023 12:01"
$ws.Range("L43").Value = 4.91
$ws.Range("M43").Value = "23/09/2023 15:28"
$ws.Range("N43").Value = 4.16
$ws.Range("O43").Value = "05/09/2023 12:01"
$ws.Range("P43").Value = 4.55
$ws.Range("Q43").Value = "23/09/2023 15:28"
$ws.Range("R43").Value = 1.83
$ws.Range("S43").Value = "05/09/2023 12:01"
$ws.Range("T43").Value = 1.66
$ws.Range("U43").Value = "23/09/2023 15:28"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/germany/bundesliga/b-monchengladbach-rb-leipzig/8M5YrEcm/"

# Row 57
$ws.Range("F57").Value = "Dortmund"
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = "Union Berlin"
$ws.Range("I57").Value = 2
$ws.Range("J57").Value = 1.52
$ws.Range("K57").Value = "23/09/2023 19:02"
$ws.Range("L57").Value = 1.71
$ws.Range("M57").Value = "07/10/2023 15:25"
$ws.Range("N57").Value = 4.52
$ws.Range("O57").Value = "23/09/2023 19:02"
$ws.Range("P57").Value = 4.17
$ws.Range("Q57").Value = "07/10/2023 15:27"
$ws.Range("R57").Value = 6.35
$ws.Range("S57").Value = "23/09/2023 19:02"
$ws.Range("T57").Value = 4.97
$ws.Range("U57").Value = "07/10/2023 15:29"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/germany/bundesliga/dortmund-union-berlin/OILWi9tI/"

# Row 59
$ws.Range("F59").Value = "Stuttgart"
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = "Wolfsburg"
$ws.Range("I59").Value = 1
$ws.Range("J59").Value = 1.96
$ws.Range("K59").Value = "23/09/2023 21:02"
$ws.Range("L59").Value = 1.88
$ws.Range("M59").Value = "07/10/2023 15:29"
$ws.Range("N59").Value = 3.8
$ws.Range("O59").Value = "23/09/2023 21:02"
$ws.Range("P59").Value = 3.97
$ws.Range("Q59").Value = "07/10/2023 15:29"
$ws.Range("R59").Value = 3.8
$ws.Range("S59").Value = "23/09/2023 21:02"
$ws.Range("T59").Value = 4.16
$ws.Range("U59").Value = "07/10/2023 15:29"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-wolfsburg/OY9f7lt5/"

# Row 60
$ws.Range("F60").Value = "Augsburg"
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = "Darmstadt"
$ws.Range("I60").Value = 2
$ws.Range("J60").Value = 1.8
$ws.Range("K60").Value = "28/09/2023 14:20"
$ws.Range("L60").Value = 1.85
$ws.Range("M60").Value = "07/10/2023 15:01"
$ws.Range("N60").Value = 3.89
$ws.Range("O60").Value = "28/09/2023 14:20"
$ws.Range("P60").Value = 3.75
$ws.Range("Q60").Value = "07/10/2023 15:29"
$ws.Range("R60").Value = 4.06
$ws.Range("S60").Value = "28/09/2023 14:20"
$ws.Range("T60").Value = 4.66
$ws.Range("U60").Value = "07/10/2023 15:28"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-darmstadt/fa5j8UQb/"

# Row 68
$ws.Range("F68").Value = "Hoffenheim"
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = "Eintracht Frankfurt"
$ws.Range("I68").Value = 3
$ws.Range("J68").Value = 2.39
$ws.Range("K68").Value = "01/10/2023 23:01"
$ws.Range("L68").Value = 2.06
$ws.Range("M68").Value = "21/10/2023 15:28"
$ws.Range("N68").Value = 3.56
$ws.Range("O68").Value = "01/10/2023 23:01"
$ws.Range("P68").Value = 3.65
$ws.Range("Q68").Value = "21/10/2023 15:26"
$ws.Range("R68").Value = 2.98
$ws.Range("S68").Value = "01/10/2023 23:01"
$ws.Range("T68").Value = 3.79
$ws.Range("U68").Value = "21/10/2023 15:28"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/germany/bundesliga/hoffenheim-eintracht-frankfurt/fDpPMTuh/"

# Row 69
$ws.Range("F69").Value = "Darmstadt"
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = "RB Leipzig"
$ws.Range("I69").Value = 3
$ws.Range("J69").Value = 6.5
$ws.Range("K69").Value = "02/10/2023 08:32"
$ws.Range("L69").Value = 6.59
$ws.Range("M69").Value = "21/10/2023 15:29"
$ws.Range("N69").Value = 5.15
$ws.Range("O69").Value = "02/10/2023 08:32"
$ws.Range("P69").Value = 5.1
$ws.Range("Q69").Value = "21/10/2023 15:29"
$ws.Range("R69").Value = 1.4
$ws.Range("S69").Value = "02/10/2023 08:32"
$ws.Range("T69").Value = 1.47
$ws.Range("U69").Value = "21/10/2023 15:28"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/germany/bundesliga/darmstadt-rb-leipzig/2wlXK7A4/"

# Row 70
$ws.Range("F70").Value = "Freiburg"
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = "Bochum"
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 1.56
$ws.Range("K70").Value = "01/10/2023 23:02"
$ws.Range("L70").Value = 1.81
$ws.Range("M70").Value = "21/10/2023 15:24"
$ws.Range("N70").Value = 4.28
$ws.Range("O70").Value = "01/10/2023 23:02"
$ws.Range("P70").Value = 3.91
$ws.Range("Q70").Value = "21/10/2023 15:27"
$ws.Range("R70").Value = 5.28
$ws.Range("S70").Value = "01/10/2023 23:02"
$ws.Range("T70").Value = 4.61
$ws.Range("U70").Value = "21/10/2023 15:26"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-bochum/rZG64nQN/"

# Row 76
$ws.Range("F76").Value = "B. Monchengladbach"
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = "Heidenheim"
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 1.93
$ws.Range("K76").Value = "10/10/2023 14:28"
$ws.Range("L76").Value = 1.69
$ws.Range("M76").Value = "28/10/2023 15:28"
$ws.Range("N76").Value = 3.9
$ws.Range("O76").Value = "10/10/2023 14:28"
$ws.Range("P76").Value = 4.44
$ws.Range("Q76").Value = "28/10/2023 15:29"
$ws.Range("R76").Value = 3.55
$ws.Range("S76").Value = "10/10/2023 14:28"
$ws.Range("T76").Value = 4.7
$ws.Range("U76").Value = "28/10/2023 15:29"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/germany/bundesliga/b-monchengladbach-heidenheim/j7lIacvd/"

# Row 77
$ws.Range("F77").Value = "Stuttgart"
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = "Hoffenheim"
$ws.Range("I77").Value = 3
$ws.Range("J77").Value = 1.77
$ws.Range("K77").Value = "10/10/2023 14:02"
$ws.Range("L77").Value = 1.67
$ws.Range("M77").Value = "28/10/2023 15:27"
$ws.Range("N77").Value = 4.11
$ws.Range("O77").Value = "10/10/2023 14:02"
$ws.Range("P77").Value = 4.43
$ws.Range("Q77").Value = "28/10/2023 15:29"
$ws.Range("R77").Value = 4.38
$ws.Range("S77").Value = "10/10/2023 14:02"
$ws.Range("T77").Value = 4.84
$ws.Range("U77").Value = "28/10/2023 15:29"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-hoffenheim/EuzDLv1F/"

# Row 78
$ws.Range("F78").Value = "Augsburg"
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = "Wolfsburg"
$ws.Range("I78").Value = 2
$ws.Range("J78").Value = 2.72
$ws.Range("K78").Value = "10/10/2023 14:02"
$ws.Range("L78").Value = 2.54
$ws.Range("M78").Value = "28/10/2023 14:56"
$ws.Range("N78").Value = 3.62
$ws.Range("O78").Value = "10/10/2023 14:02"
$ws.Range("P78").Value = 3.62
$ws.Range("Q78").Value = "28/10/2023 14:56"
$ws.Range("R78").Value = 2.44
$ws.Range("S78").Value = "10/10/2023 14:02"
$ws.Range("T78").Value = 2.8
$ws.Range("U78").Value = "28/10/2023 14:56"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-wolfsburg/d4u8MKo9/"

# Row 79
$ws.Range("F79").Value = "Bayern Munich"
$ws.Range("G79").Value = 8
$ws.Range("H79").Value = "Darmstadt"
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 1.07
$ws.Range("K79").Value = "10/10/2023 14:25"
$ws.Range("L79").Value = 1.06
$ws.Range("M79").Value = "28/10/2023 15:16"
$ws.Range("N79").Value = 12.93
$ws.Range("O79").Value = "10/10/2023 14:25"
$ws.Range("P79").Value = 17.61
$ws.Range("Q79").Value = "28/10/2023 15:29"
$ws.Range("R79").Value = 19.15
$ws.Range("S79").Value = "10/10/2023 14:25"
$ws.Range("T79").Value = 32.15
$ws.Range("U79").Value = "28/10/2023 15:29"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayern-munich-darmstadt/p2wxKRPA/"

# Row 84
$ws.Range("F84").Value = "FC Koln"
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = "Augsburg"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 1.9
$ws.Range("K84").Value = "22/10/2023 12:02"
$ws.Range("L84").Value = 2.05
$ws.Range("M84").Value = "04/11/2023 15:29"
$ws.Range("N84").Value = 3.73
$ws.Range("O84").Value = "22/10/2023 12:02"
$ws.Range("P84").Value = 3.89
$ws.Range("Q84").Value = "04/11/2023 15:29"
$ws.Range("R84").Value = 3.8
$ws.Range("S84").Value = "22/10/2023 12:02"
$ws.Range("T84").Value = 3.52
$ws.Range("U84").Value = "04/11/2023 15:29"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-augsburg/CzIaCzO1/"

# Row 85
$ws.Range("F85").Value = "Freiburg"
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = "B. Monchengladbach"
$ws.Range("I85").Value = 3
$ws.Range("J85").Value = 1.77
$ws.Range("K85").Value = "22/10/2023 12:02"
$ws.Range("L85").Value = 2.03
$ws.Range("M85").Value = "04/11/2023 15:19"
$ws.Range("N85").Value = 3.99
$ws.Range("O85").Value = "22/10/2023 12:02"
$ws.Range("P85").Value = 3.85
$ws.Range("Q85").Value = "04/11/2023 15:19"
$ws.Range("R85").Value = 4.55
$ws.Range("S85").Value = "22/10/2023 12:02"
$ws.Range("T85").Value = 3.64
$ws.Range("U85").Value = "04/11/2023 15:19"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-b-monchengladbach/MoPnFxvq/"

# Row 86
$ws.Range("F86").Value = "Hoffenheim"
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = "Bayer Leverkusen"
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 3.85
$ws.Range("K86").Value = "22/10/2023 12:02"
$ws.Range("L86").Value = 5.53
$ws.Range("M86").Value = "04/11/2023 15:29"
$ws.Range("N86").Value = 4.05
$ws.Range("O86").Value = "22/10/2023 12:02"
$ws.Range("P86").Value = 5.01
$ws.Range("Q86").Value = "04/11/2023 15:27"
$ws.Range("R86").Value = 1.81
$ws.Range("S86").Value = "22/10/2023 12:02"
$ws.Range("T86").Value = 1.53
$ws.Range("U86").Value = "04/11/2023 15:29"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/germany/bundesliga/hoffenheim-bayer-leverkusen/jRG3Bfw8/"

# Row 87
$ws.Range("F87").Value = "Mainz"
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = "RB Leipzig"
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 3.62
$ws.Range("K87").Value = "22/10/2023 12:02"
$ws.Range("L87").Value = 4.57
$ws.Range("M87").Value = "04/11/2023 15:28"
$ws.Range("N87").Value = 3.78
$ws.Range("O87").Value = "22/10/2023 12:02"
$ws.Range("P87").Value = 4.04
$ws.Range("Q87").Value = "04/11/2023 15:28"
$ws.Range("R87").Value = 1.93
$ws.Range("S87").Value = "22/10/2023 12:02"
$ws.Range("T87").Value = 1.78
$ws.Range("U87").Value = "04/11/2023 15:28"
$ws.Range("V87").Value = "https://www.betexplorer.com/football/germany/bundesliga/mainz-rb-leipzig/r1DeDG8e/"

# --- 2) Append two new rows at the end (90 and 91) ---
# Copy formatting from the last existing data row (89) for the styled columns A and E
$ws.Range("A89").Copy() | Out-Null
$ws.Range("A90").PasteSpecial(-4122) | Out-Null
$ws.Range("A89").Copy() | Out-Null
$ws.Range("A91").PasteSpecial(-4122) | Out-Null
$ws.Range("E89").Copy() | Out-Null
$ws.Range("E90").PasteSpecial(-4122) | Out-Null
$ws.Range("E89").Copy() | Out-Null
$ws.Range("E91").PasteSpecial(-4122) | Out-Null

# Row 90
$ws.Range("A90").Value = 89
$ws.Range("B90").Value = "germany"
$ws.Range("C90").Value = "bundesliga"
$ws.Range("D90").Value = "2023-2024"
$ws.Range("E90").Value = 45235.64583333334
$ws.Range("F90").Value = "Wolfsburg"
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = "Werder Bremen"
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 1.61
$ws.Range("K90").Value = "22/10/2023 12:02"
$ws.Range("L90").Value = 1.75
$ws.Range("M90").Value = "05/11/2023 15:29"
$ws.Range("N90").Value = 4.44
$ws.Range("O90").Value = "22/10/2023 12:02"
$ws.Range("P90").Value = 4.31
$ws.Range("Q90").Value = "05/11/2023 15:29"
$ws.Range("R90").Value = 5.26
$ws.Range("S90").Value = "22/10/2023 12:02"
$ws.Range("T90").Value = 4.39
$ws.Range("U90").Value = "05/11/2023 15:29"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/germany/bundesliga/wolfsburg-werder-bremen/4AEiEdgk/"

# Row 91
$ws.Range("A91").Value = 90
$ws.Range("B91").Value = "germany"
$ws.Range("C91").Value = "bundesliga"
$ws.Range("D91").Value = "2023-2024"
$ws.Range("E91").Value = 45235.72916666666
$ws.Range("F91").Value = "Heidenheim"
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = "Stuttgart"
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3.38
$ws.Range("K91").Value = "22/10/2023 15:56"
$ws.Range("L91").Value = 5.3
$ws.Range("M91").Value = "05/11/2023 17:28"
$ws.Range("N91").Value = 3.8
$ws.Range("O91").Value = "22/10/2023 15:56"
$ws.Range("P91").Value = 4.25
$ws.Range("Q91").Value = "05/11/2023 17:28"
$ws.Range("R91").Value = 2.01
$ws.Range("S91").Value = "22/10/2023 15:56"
$ws.Range("T91").Value = 1.65
$ws.Range("U91").Value = "05/11/2023 17:28"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/germany/bundesliga/heidenheim-vfb-stuttgart/faspKns8/"

